$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Survey")

# Row 1 - header. B1 ("Question") is entered before A1 ("ID") so the shared
# string table is built in the same first-use order as the target workbook.
$ws.Cells.Item(1,2).Value = "Question"
$ws.Cells.Item(1,1).Value = "ID"
$ws.Cells.Item(1,3).Value = 1
$ws.Cells.Item(1,4).Value = 2
$ws.Cells.Item(1,5).Value = 3
$ws.Cells.Item(1,6).Value = 4
$ws.Cells.Item(1,7).Value = 5
$ws.Cells.Item(1,8).Value = "Solution"

# Row 2 - Pizza question
$ws.Cells.Item(2,1).Value = 1
$ws.Cells.Item(2,2).Value = "Pizza?"
$ws.Cells.Item(2,3).Value = "Tomaten"
$ws.Cells.Item(2,4).Value = "Schinken"
$ws.Cells.Item(2,5).Value = "Käse"
$ws.Cells.Item(2,6).Value = "Salami"
$ws.Cells.Item(2,7).Value = "Ananas"
$ws.Cells.Item(2,8).Value = 3

# Row 3 - Döner question
$ws.Cells.Item(3,1).Value = 2
$ws.Cells.Item(3,2).Value = "Döner?"
$ws.Cells.Item(3,3).Value = "Scharf"
$ws.Cells.Item(3,4).Value = "Käse"
$ws.Cells.Item(3,5).Value = "Vegetarisch"
$ws.Cells.Item(3,6).Value = "Fleisch"
$ws.Cells.Item(3,8).Value = "1;2"

# Row 4 - Nudeln question
$ws.Cells.Item(4,1).Value = 3
$ws.Cells.Item(4,2).Value = "Nudeln?"
$ws.Cells.Item(4,3).Value = "Tomaten"
$ws.Cells.Item(4,4).Value = "Sahne"
$ws.Cells.Item(4,8).Value = 1

# Column widths: column A narrow (ID numbers), columns B:H wide (answers)
$ws.Columns.Item(1).ColumnWidth = 2
$ws.Range("B1:H1").EntireColumn.ColumnWidth = 11.83

# Selection / active sheet: Survey becomes the active tab with A6 selected
$ws.Range("A6").Select()
